# build_4_future_2015.16x9.pptx edit script
# - refresh the "Update automatically" date placeholder (5/10/15 -> 5/19/15)
#   on the slide master, every slide layout, and the notes master
# - fix the speed-of-light math error on slide 14 (299792.5 -> 299793,
#   455684.5 -> 197231)
# - collapse the split "Boyd " / "Multerer" runs on slide 2 into one run

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholders: "Update automatically" field text 5/10/15 -> 5/19/15
# ---------------------------------------------------------------------------

function Set-DatePlaceholderText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "5/10/15") {
                $tr.Text = "5/19/15"
            }
        }
    }
}

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes

# Every slide layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes
}

# Notes master
if ($p.HasNotesMaster) {
    Set-DatePlaceholderText $p.NotesMaster.Shapes
} else {
    Set-DatePlaceholderText $p.NotesMaster.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 14: fix "sill" math error in speed-of-light figures
# ---------------------------------------------------------------------------

$s14 = $p.Slides.Item(14)
$contentShape = $s14.Shapes.Item(2)
$tr14 = $contentShape.TextFrame.TextRange

# Paragraph 5: "299792.5 km/s" -> "299793 km/s"
$para5 = $tr14.Paragraphs(5, 1)
$numPart5 = $para5.Characters(1, 9)
$numPart5.Text = "299793 "

# Paragraph 7: "455684.5 km/s" -> "197231 km/s"
$para7 = $tr14.Paragraphs(7, 1)
$numPart7 = $para7.Characters(1, 9)
$numPart7.Text = "197231 "

# ---------------------------------------------------------------------------
# 3. Slide 2: merge split "Boyd " + "Multerer" runs into "Boyd Multerer"
# ---------------------------------------------------------------------------

$s2 = $p.Slides.Item(2)
$subtitle = $s2.Shapes.Item(2)
$trSub = $subtitle.TextFrame.TextRange
$paraName = $trSub.Paragraphs(1, 1)

# Force the engine to rebuild this as a single run: an identity assignment
# is a no-op, so first swap in an unrelated placeholder value, then the
# final text, so the run gets rewritten as one contiguous run.
$paraName.Text = "ZZZZZZZZZZZZZ"
$paraName.Text = "Boyd Multerer"
